$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Clear out the old "Good Drivers" tail rows (31-49) entirely ----
$ws.Range("A31:J49").Clear()

# ---- Bad Drivers table (rows 3-6) ----
$ws.Range("A3").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.27.1"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 782
$ws.Range("D3").Value = 94.40000000000001

$ws.Range("A4").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.30.1"
$ws.Range("B4").Value = 24
$ws.Range("C4").Value = 1346
$ws.Range("D4").Value = 98.8

$ws.Range("A5").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.18.2"
$ws.Range("B5").Value = 21
$ws.Range("C5").Value = 1911
$ws.Range("D5").Value = 98.8

$ws.Range("A6").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.80.0.7"
$ws.Range("B6").Value = 60
$ws.Range("C6").Value = 2598
$ws.Range("D6").Value = 98.90000000000001

# ---- Totals row (row 7) ----
$ws.Range("B7").Value = 107
$ws.Range("C7").Value = 6637

# ---- Good Drivers table (rows 15-30) ----
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B15").Value = 56018
$ws.Range("D15").Value = 100
$ws.Range("E15").ClearContents()

$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B16").Value = 34244
$ws.Range("D16").Value = 100
$ws.Range("E16").ClearContents()

$ws.Range("A17").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.11.3"
$ws.Range("B17").Value = 161874
$ws.Range("D17").Value = 100
$ws.Range("E17").ClearContents()

$ws.Range("A18").Value = "Intel(R) Dual Band Wireless-AC 8265 - 22.0.1.1"
$ws.Range("B18").Value = 52096
$ws.Range("D18").Value = 100
$ws.Range("E18").ClearContents()

$ws.Range("A19").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B19").Value = 442178
$ws.Range("D19").Value = 99.90000000000001
$ws.Range("E19").Value = "'" + "2024-11-10"

$ws.Range("A20").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B20").Value = 77849
$ws.Range("D20").Value = 99.90000000000001
$ws.Range("E20").Value = "'" + "2021-08-18"

$ws.Range("A21").Value = "Intel(R) Dual Band Wireless-AC 8265 - 22.30.0.11"
$ws.Range("B21").Value = 170510
$ws.Range("D21").Value = 99.90000000000001
$ws.Range("E21").Value = "'" + "2021-01-19"

$ws.Range("A22").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.12.5"
$ws.Range("B22").Value = 141909
$ws.Range("D22").Value = 99.90000000000001
$ws.Range("E22").Value = "'" + "2021-01-19"

$ws.Range("A23").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B23").Value = 59673
$ws.Range("D23").Value = 100
$ws.Range("E23").Value = "'" + "2020-08-05"

$ws.Range("A24").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B24").Value = 113652
$ws.Range("D24").Value = 100
$ws.Range("E24").Value = "'" + "2019-12-14"

$ws.Range("A25").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.8.1"
$ws.Range("B25").Value = 48540
$ws.Range("D25").Value = 100
$ws.Range("E25").Value = "'" + "2019-09-05"

$ws.Range("A26").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.5.2"
$ws.Range("B26").Value = 184564
$ws.Range("D26").Value = 99.90000000000001
$ws.Range("E26").Value = "'" + "2019-08-25"

$ws.Range("A27").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.10.2"
$ws.Range("B27").Value = 20227
$ws.Range("D27").Value = 100
$ws.Range("E27").Value = "'" + "2019-05-11"

$ws.Range("A28").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.9.1"
$ws.Range("B28").Value = 34065
$ws.Range("D28").Value = 100
$ws.Range("E28").Value = "'" + "2019-04-28"

$ws.Range("A29").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.50.0.4"
$ws.Range("B29").Value = 14221
$ws.Range("D29").Value = 100
$ws.Range("E29").Value = "'" + "2018-05-08"

$ws.Range("A30").Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.30.1.2"
$ws.Range("B30").Value = 23765
$ws.Range("D30").Value = 100
$ws.Range("E30").Value = "'" + "2018-01-09"

# ---- Force dimension to extend through column J / row 35 (template spacer columns) ----
$ws.Range("J35").NumberFormat = "General"
